$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36 (4/2/2020) now has an actual reported value instead of a forecast.
# Convert I36 from the forecast formula to a hardcoded "actual" value and
# give it the "actual data" fill style (same style used by I30:I35).
$ws.Range("I36").Value = 244877
$ws.Range("I36").Style = $ws.Range("I35").Style

# I37's shared forecast formula now anchors off I36 (shifts window by one row).
$ws.Range("I37").Formula = "=I36*(1+AVERAGE(M34:M36))"

# Column H got a bit narrower.
$ws.Range("H1").EntireColumn.ColumnWidth = 8.6640625

# Selection moved to I37.
$ws.Range("I37").Select()
